$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header text (December 2021 running total label)
$ws.Range("B1").Value = "December 2021 (through December 22)"

# Increment existing December counts (new carjacking report added to multiple
# December columns across years) per the source diff
$ws.Range("B2").Value = 2   # was 1
$ws.Range("B3").Value = 8   # was 7
$ws.Range("Z4").Value = 6   # was 5
$ws.Range("N6").Value = 10   # was 9
$ws.Range("Z6").Value = 4   # was 3
$ws.Range("AX6").Value = 8   # was 7
$ws.Range("AL7").Value = 11   # was 10
$ws.Range("N8").Value = 3   # was 2
$ws.Range("N9").Value = 3   # was 2
$ws.Range("AL9").Value = 2   # was 1
$ws.Range("AX12").Value = 3   # was 2
$ws.Range("B13").Value = 3   # was 2
$ws.Range("BJ13").Value = 7   # was 6
$ws.Range("AX15").Value = 2   # was 1
$ws.Range("B18").Value = 5   # was 4
$ws.Range("Z19").Value = 2   # was 1
$ws.Range("N23").Value = 3   # was 2
$ws.Range("B39").Value = 2   # was 1
$ws.Range("Z40").Value = 2   # was 1
$ws.Range("AL65").Value = 2   # was 1

# New cells introduced where no prior December count existed
$ws.Range("BV18").Value = 1
$ws.Range("N33").Value = 1
$ws.Range("BJ34").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("BV51").Value = 1
$ws.Range("AX53").Value = 1
$ws.Range("AX55").Value = 1
$ws.Range("Z66").Value = 1
$ws.Range("BJ71").Value = 1
$ws.Range("N95").Value = 1
$ws.Range("AX97").Value = 1

# Finally, rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-12-22"

